$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '276.87'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.82%'
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.26%'
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.855'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.05%'
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.19%'
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.939'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.54%'
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.206'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-6.60%'
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8766'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.96%'
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '4.50%'
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05124'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '3.30%'
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.81%'
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02961'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.03%'
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.08981'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.54%'
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001571'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.41%'
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006369'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.91%'
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006097'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.80%'
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.476'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.83%'
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.306'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.60%'
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.275'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.37%'
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.30%'
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.69%'
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.903'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-0.09%'
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04426'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.72%'
$ws.Range("E23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001176'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.26%'
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.003862'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-9.24%'
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001300'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '8.28%'
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '15.08%'
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04155'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.47%'
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006802'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.90%'
$ws.Range("E41").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001959'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-6.73%'
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01191'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '11.08%'
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005298'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-0.08%'
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '13.27%'
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01852'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-7.45%'
$ws.Range("E47").Style = "Normal"
